$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content (old Week/Task table)
$ws.Cells.Clear() | Out-Null

# Column A needs to be wider for the new task descriptions
$ws.Columns.Item(1).ColumnWidth = 34.7271205357

# Populate the new task/hours table
$ws.Range("A7").Value = 'Tasks '
$ws.Range("B7").Value = 'Hours'

$ws.Range("A8").Value = 'Task 1 Week 1'
$ws.Range("A8").Font.Bold = $true

$ws.Range("A9").Value = 'Set up a Git Repo'
$ws.Range("B9").Value = '10 min'

$ws.Range("A10").Value = 'Set up Trello for each task'
$ws.Range("B10").Value = '20 min'

$ws.Range("A11").Value = 'Create a sample WPF program'
$ws.Range("B11").Value = '30 min'

$ws.Range("A12").Value = 'Follow online resoure for tips'
$ws.Range("B12").Value = '10 min'

$ws.Range("A13").Value = 'Task 2 Week 1'
$ws.Range("A13").Font.Bold = $true

$ws.Range("A14").Value = 'Create the methods for URI conn. '
$ws.Range("B14").Value = ' 5 hrs'

$ws.Range("A15").Value = 'Figure out how to use oauth 2.0'
$ws.Range("B15").Value = ' 5 hrs'

$ws.Range("A16").Value = 'Connect parser and library to sample console'
$ws.Range("B16").Value = ' 3 hrs'

$ws.Range("A17").Value = 'Fix bugs that shown up'
$ws.Range("B17").Value = ' 3  hrs'

$ws.Range("A18").Value = 'Task 3 : Week 2'
$ws.Range("A18").Font.Bold = $true

$ws.Range("A19").Value = 'Connect parser and library to WPF'
$ws.Range("B19").Value = '5 hrs'

$ws.Range("A20").Value = 'Properly implement library'
$ws.Range("B20").Value = '5 hrs'

$ws.Range("A21").Value = 'Fix bugs in the from connecting library and  WPf'
$ws.Range("B21").Value = '10 hrs'

$ws.Range("A22").Value = 'Task 4: Week 3'
$ws.Range("A22").Font.Bold = $true

$ws.Range("A23").Value = 'Write up documentation'
$ws.Range("B23").Value = '1 hr'

$ws.Range("A24").Value = 'Write up log for testers'
$ws.Range("B24").Value = '1 hr'

# Match the selection left by the author
$ws.Range("B24").Select() | Out-Null

# Page orientation switched to portrait
$ws.PageSetup.Orientation = 1
